$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "host_organization" column (column H). Excel shifts every
# column to its right one position to the left, so "issn_l" (was I) is
# now H, "url" (was J) is now I, and so on through to "is_retracted"
# (was AH, now AG).
$ws.Range("H1").EntireColumn.Delete()

# The source data used to populate this row was fixed, so a few of the
# shifted values need correcting to their real (previously "N/A" /
# stale) values.
$ws.Range("F2").Value2 = "IOP Conference Series: Materials Science and Engineering"
$ws.Range("G2").Value2 = "https://openalex.org/S4210189194"
$ws.Range("H2").Value2 = "1757-8981"
# Leading apostrophe forces this to be stored as text "FALSE" rather
# than the Boolean FALSE (matches the other TRUE/FALSE text cells on
# this sheet, which are all plain strings, not typed booleans). Reset
# the style afterwards so the quote-prefix marker doesn't change the
# cell's formatting away from the sheet's default (unstyled) look.
$ws.Range("U2").Value2 = "'FALSE"
$ws.Range("U2").Style = "Normal"
